$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G. This shifts the existing
# @KingEffectIcon / KingEffectText / FavourText columns from G:I to H:J.
$ws.Columns("G").EntireColumn.Insert()

# New "VictoryText" column header and King's victory condition text.
$ws.Range("G1").Value = "VictoryText"
$ws.Range("G2").Value = "Wins if alive."

# Pick up the header/body formatting used by the neighbouring columns
# (the raw column insert otherwise clones column A's formatting).
$ws.Range("H1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H2:H4").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the column width used by the neighbouring King-effect columns,
# and give the shifted FavourText column (now J) its own slightly
# narrower width, same as in the finished template.
$ws.Columns("G").ColumnWidth = 31.95
$ws.Columns("J").ColumnWidth = 31.45

# Selection left on G3 by whoever was editing this template.
$ws.Range("G3").Select()

Write-Output "done"
